# activity-rule.xlsx: rename taskKey "e.init" to "e.logout" for the
# "领导审批通过" / logout related rule rows (C21:C26 on sheet DATA-CDATA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 21; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -eq "e.init") {
        $cell.Value = "e.logout"
    }
}
